# "Analisis y Diseno terminado"
# Fix a typo in the header of column G (cantiadad_estudiantes -> cantidad_estudiantes)
# and append three new course rows (introduccion, algebra, quimica) to the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header typo ---
$ws.Cells.Item(1, 7).Value = "cantidad_estudiantes"

# --- Row 12: introduccion ---
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 110
$ws.Cells.Item(12, 3).Value = "introduccion"
$ws.Cells.Item(12, 4).Value = 12
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = "miercoles"
$ws.Cells.Item(12, 7).Value = 17

# --- Row 13: algebra ---
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 110
$ws.Cells.Item(13, 6).Value = "sabado"
$ws.Cells.Item(13, 3).Value = "algebra"
$ws.Cells.Item(13, 4).Value = 10
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 7).Value = 14

# --- Row 14: quimica ---
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = 110
$ws.Cells.Item(14, 3).Value = "quimica"
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = "viernes"
$ws.Cells.Item(14, 7).Value = 23

# --- Cursor / selection state, as left by the author ---
[void]$ws.Cells.Item(8, 9).Select()

# --- Page orientation set to portrait for printing ---
$ws.PageSetup.Orientation = 1
